# Apply the update described by the diff:
#  - Insert 3 new product rows into the table (keeping alphabetical Arabic order)
#  - Renumber the index column, refresh the totals row and the timestamp footer
#
# Excel constants used (no enum available, so literal values are used):
#   xlShiftDown        = -4121
#   xlPasteFormats      = -4122
#   xlPasteValues       = -4163

$xlShiftDown   = -4121
$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert the three new blank rows (bottom-most target first shifts less
#    work, but since each insert changes the position of subsequent targets,
#    we insert from the top down and recompute the target row every time).
# ---------------------------------------------------------------------------

# New row for "سرنجات 10 سم" goes right before the current "سرنجات 5 سم" row (row 11)
$ws.Range("A11:N11").Insert($xlShiftDown)
$ws.Range("A10:N10").Copy()
$ws.Range("A11:N11").PasteSpecial($xlPasteFormats)

# New row for "شفرات فينوس حريمي " goes right before "كالونا " (now row 13)
$ws.Range("A13:N13").Insert($xlShiftDown)
$ws.Range("A12:N12").Copy()
$ws.Range("A13:N13").PasteSpecial($xlPasteFormats)

# New row for "كريم فيرند لافلي الصغير" goes right after "كالونا " (now row 15)
$ws.Range("A15:N15").Insert($xlShiftDown)
$ws.Range("A14:N14").Copy()
$ws.Range("A15:N15").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# Row heights for rows 11-16 (per target layout; Excel recalculates/auto-fits
# these based on wrapped text, so every row height is restated explicitly)
$ws.Rows.Item(11).RowHeight = 24.75
$ws.Rows.Item(12).RowHeight = 25.5
$ws.Rows.Item(13).RowHeight = 25.5
$ws.Rows.Item(14).RowHeight = 24.75
$ws.Rows.Item(15).RowHeight = 25.5
$ws.Rows.Item(16).RowHeight = 24.75

# ---------------------------------------------------------------------------
# 2. Re-create the merged cells for the 3 new rows (lost on insert)
# ---------------------------------------------------------------------------
$ws.Range("B11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()

$ws.Range("B13:G13").Merge()
$ws.Range("H13:K13").Merge()
$ws.Range("L13:M13").Merge()

$ws.Range("B15:G15").Merge()
$ws.Range("H15:K15").Merge()
$ws.Range("L15:M15").Merge()

# ---------------------------------------------------------------------------
# 3. Fill in the values for the new rows and refresh the renumbered index
#    column plus the totals that changed because of the new products.
# ---------------------------------------------------------------------------

# Row 11 : سرنجات 10 سم
$ws.Cells.Item(11,1).Value  = 8
$ws.Cells.Item(11,2).Value  = "سرنجات 10 سم"
$ws.Cells.Item(11,8).Value  = "-2:0"
$ws.Cells.Item(11,12).Value = 8
$ws.Cells.Item(11,14).Value = "2:0"

# Row 12 : سرنجات 5 سم (already existed, only index + totals move)
$ws.Cells.Item(12,1).Value  = 9

# Row 13 : شفرات فينوس حريمي
$ws.Cells.Item(13,1).Value  = 10
$ws.Cells.Item(13,2).Value  = "شفرات فينوس حريمي "
$ws.Cells.Item(13,8).Value  = "16:0"
$ws.Cells.Item(13,12).Value = 40
$ws.Cells.Item(13,14).Value = "2:0"

# Row 14 : كالونا (already existed, only index moves)
$ws.Cells.Item(14,1).Value  = 11

# Row 15 : كريم فيرند لافلي الصغير
$ws.Cells.Item(15,1).Value  = 12
$ws.Cells.Item(15,2).Value  = "كريم فيرند لافلي الصغير"
$ws.Cells.Item(15,8).Value  = "6:0"
$ws.Cells.Item(15,12).Value = 20
$ws.Cells.Item(15,14).Value = "1:0"

# Row 16 : محلول ملح (already existed, only index moves)
$ws.Cells.Item(16,1).Value  = 13

# ---------------------------------------------------------------------------
# 4. Totals row and footer (now shifted from rows 14/15 to rows 17/18)
# ---------------------------------------------------------------------------
$ws.Cells.Item(17,11).Value = 446
$ws.Rows.Item(17).RowHeight = 26.25

Write-Output "edit applied"
